{"js": "// Mute button and Welcome Screen Login\n//\n// Net content change: two obsolete objective rows are removed from their\n// respective tables:\n//   1) \"Submit the project on time\" / \"3\"  (Quality of Work table)\n//   2) \"Find a way to show the number of songs in our playlist\" / \"0\"\n//      (Objectives table for the \"-Above & Beyond\" section)\n//\n// We locate the rows by their first-cell text (robust against table\n// position) and delete the whole row via TableRow.delete().\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Submit the project on time\",\n  \"Find a way to show the number of songs in our playlist\",\n]);\n\nconst rowsToDelete = [];\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < rows.items.length; j++) {\n    const row = rows.items[j];\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    if (cells.items.length === 0) continue;\n\n    const firstCell = cells.items[0];\n    firstCell.load(\"value\");\n    await context.sync();\n\n    const text = (firstCell.value || \"\").trim();\n    if (targets.has(text)) {\n      rowsToDelete.push(row);\n    }\n  }\n}\n\nfor (const row of rowsToDelete) {\n  row.delete();\n}\n\nawait context.sync();\n", "ps1": "# Mute button and Welcome Screen Login\n#\n# Net content change: two obsolete objective rows are removed from their\n# respective tables:\n#   1) \"Submit the project on time\" / \"3\"  (Quality of Work table)\n#   2) \"Find a way to show the number of songs in our playlist\" / \"0\"\n#      (Objectives table for the \"-Above & Beyond\" section)\n#\n# We scan every table in the document, find rows whose first cell matches\n# one of the target strings, and delete those rows. We collect matches\n# first and delete from the highest row index down so earlier indices\n# stay valid while deleting.\n\n$d = $word.ActiveDocument\n\n$targets = @(\n  \"Submit the project on time\",\n  \"Find a way to show the number of songs in our playlist\"\n)\n\nfor ($ti = 1; $ti -le $d.Tables.Count; $ti++) {\n  $tbl = $d.Tables.Item($ti)\n  $rowsToDelete = @()\n\n  for ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $cellText = $tbl.Cell($r, 1).Range.Text\n    $cleanText = $cellText.TrimEnd([char]7).TrimEnd([char]13).Trim()\n    if ($targets -contains $cleanText) {\n      $rowsToDelete += $r\n    }\n  }\n\n  for ($i = $rowsToDelete.Count - 1; $i -ge 0; $i--) {\n    $tbl.Rows.Item($rowsToDelete[$i]).Delete()\n  }\n}\n"}
